# Fix a bug in featureBar: the rows in the data table were shuffled in the
# wrong order. Re-apply the correct row ordering for rows 3-21 (A:F),
# matching the intended feature-bar bucket ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3  = @(401, 9, 48, 67, 75, 45)
    4  = @(601, 9, 60, 67, 60, 42)
    5  = @(1203, 3, 15, 15, 15, 15)
    6  = @(902, 1, 0, 0, 0, 0)
    7  = @(701, 3, 90, 45, 97, 15)
    8  = @(201, 9, 30, 15, 45, 30)
    9  = @(801, 3, 67, 65, 52, 45)
    10 = @(1202, 2, 10, 10, 10, 10)
    11 = @(1001, 18, 30, 75, 60, 72)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(1201, 2, 10, 10, 10, 10)
    14 = @(901, 16, 15, 45, 60, 60)
    15 = @(301, 6, 45, 30, 60, 45)
    16 = @(1101, 0, 15, 30, 30, 0)
    17 = @(2, 0, 2, 2, 2, 2)
    19 = @(1, 0, 2, 2, 2, 2)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(3, 0, 3, 3, 3, 3)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
